# Update the "想去人数" (interest count) column F across the four sheets.
# Each entry increments the previously published count, matching the
# gh-pages data refresh captured in the commit (output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

function Set-Count {
    param(
        [string]$SheetName,
        [string]$CellRef,
        [double]$OldValue,
        [double]$NewValue
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $cell = $ws.Range($CellRef)

    $current = $cell.Value2
    if ($current -ne $OldValue) {
        Write-Host ("WARNING: " + $SheetName + "!" + $CellRef + " was " + $current + ", expected " + $OldValue)
    }

    $cell.Value = $NewValue
}

# 展览 (Exhibitions)
Set-Count "展览" "F12" 1818 1819
Set-Count "展览" "F18" 528  529
Set-Count "展览" "F26" 1047 1048
Set-Count "展览" "F27" 4617 4618
Set-Count "展览" "F31" 0    1
Set-Count "展览" "F32" 176  177

# 演出 (Performances)
Set-Count "演出" "F32" 485 486

# 本地生活 (Local life)
Set-Count "本地生活" "F11" 901 903
Set-Count "本地生活" "F13" 42  44
Set-Count "本地生活" "F14" 65  67
Set-Count "本地生活" "F16" 322 323

# 全部类型 (All types - combined view)
Set-Count "全部类型" "F8"  901  903
Set-Count "全部类型" "F10" 42   44
Set-Count "全部类型" "F11" 42   44
Set-Count "全部类型" "F12" 65   67
Set-Count "全部类型" "F13" 65   67
Set-Count "全部类型" "F18" 1818 1819
Set-Count "全部类型" "F35" 1047 1048
Set-Count "全部类型" "F37" 322  323
Set-Count "全部类型" "F41" 485  486
Set-Count "全部类型" "F46" 176  177
